$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 23.20031648293105
$ws.Cells.Item(2, 3).Value = 9.121142890126858
$ws.Cells.Item(2, 4).Value = 8.180230647678927
$ws.Cells.Item(2, 5).Value = 9.726312213286635
$ws.Cells.Item(2, 6).Value = 42.38014468077433
$ws.Cells.Item(2, 12).Value = 10.02731762449565
$ws.Cells.Item(2, 13).Value = 18.44414851031817
$ws.Cells.Item(2, 14).Value = 21.31830534023515
$ws.Cells.Item(3, 2).Value = 22.79640028415665
$ws.Cells.Item(3, 3).Value = 8.520819702750796
$ws.Cells.Item(3, 4).Value = 8.20335936488641
$ws.Cells.Item(3, 5).Value = 9.708442519933733
$ws.Cells.Item(3, 6).Value = 42.07396438172254
$ws.Cells.Item(3, 12).Value = 10.04054932860222
$ws.Cells.Item(3, 13).Value = 18.37493880585774
$ws.Cells.Item(3, 14).Value = 21.37021744265667
$ws.Cells.Item(4, 2).Value = 22.55225410524609
$ws.Cells.Item(4, 3).Value = 8.129543559500849
$ws.Cells.Item(4, 4).Value = 8.218590814447213
$ws.Cells.Item(4, 5).Value = 9.697215147095093
$ws.Cells.Item(4, 6).Value = 41.89694252176102
$ws.Cells.Item(4, 12).Value = 10.05016902832889
$ws.Cells.Item(4, 13).Value = 18.33678581621582
$ws.Cells.Item(4, 14).Value = 21.40404577617458
$ws.Cells.Item(5, 2).Value = 22.45387329883227
$ws.Cells.Item(5, 3).Value = 7.964377621417144
$ws.Cells.Item(5, 4).Value = 8.225056415159044
$ws.Cells.Item(5, 5).Value = 9.692575405310716
$ws.Cells.Item(5, 6).Value = 41.82761757209413
$ws.Cells.Item(5, 12).Value = 10.05446503871258
$ws.Cells.Item(5, 13).Value = 18.3223397775157
$ws.Cells.Item(5, 14).Value = 21.4183221187192
$ws.Cells.Item(6, 2).Value = 22.43760839616037
$ws.Cells.Item(6, 3).Value = 7.936606157332932
$ws.Cells.Item(6, 4).Value = 8.226145633696444
$ws.Cells.Item(6, 5).Value = 9.691801075610893
$ws.Cells.Item(6, 6).Value = 41.81627755937322
$ws.Cells.Item(6, 12).Value = 10.05520108717847
$ws.Cells.Item(6, 13).Value = 18.32000781733118
$ws.Cells.Item(6, 14).Value = 21.42072234125611
$ws.Cells.Item(7, 2).Value = 22.55092262584575
$ws.Cells.Item(7, 3).Value = 8.127339245173038
$ws.Cells.Item(7, 4).Value = 8.218676965071349
$ws.Cells.Item(7, 5).Value = 9.697152835420161
$ws.Cells.Item(7, 6).Value = 41.89599612586279
$ws.Cells.Item(7, 12).Value = 10.05022544400861
$ws.Cells.Item(7, 13).Value = 18.33658651967619
$ws.Cells.Item(7, 14).Value = 21.40423632434943
$ws.Cells.Item(8, 2).Value = 23.06032937887691
$ws.Cells.Item(8, 3).Value = 8.918847913960878
$ws.Cells.Item(8, 4).Value = 8.187991328454093
$ws.Cells.Item(8, 5).Value = 9.720203551007844
$ws.Cells.Item(8, 6).Value = 42.27232972341336
$ws.Cells.Item(8, 12).Value = 10.03156949000334
$ws.Cells.Item(8, 13).Value = 18.41939007785535
$ws.Cells.Item(8, 14).Value = 21.33579884677151
$ws.Cells.Item(9, 2).Value = 24.08372626987857
$ws.Cells.Item(9, 3).Value = 10.29183033728997
$ws.Cells.Item(9, 4).Value = 8.136011402257703
$ws.Cells.Item(9, 5).Value = 9.763394810530185
$ws.Cells.Item(9, 6).Value = 43.09479510621729
$ws.Cells.Item(9, 12).Value = 10.00685491608024
$ws.Cells.Item(9, 13).Value = 18.61573043672768
$ws.Cells.Item(9, 14).Value = 21.21711358252183
$ws.Cells.Item(10, 2).Value = 24.84233232133102
$ws.Cells.Item(10, 3).Value = 11.19242253867698
$ws.Cells.Item(10, 4).Value = 8.102844657830772
$ws.Cells.Item(10, 5).Value = 9.793923492208648
$ws.Cells.Item(10, 6).Value = 43.74674981790653
$ws.Cells.Item(10, 12).Value = 9.995938068887705
$ws.Cells.Item(10, 13).Value = 18.77994718712095
$ws.Cells.Item(10, 14).Value = 21.13939850484969
$ws.Cells.Item(11, 2).Value = 25.18721969605706
$ws.Cells.Item(11, 3).Value = 11.57888310752609
$ws.Cells.Item(11, 4).Value = 8.088853209888113
$ws.Cells.Item(11, 5).Value = 9.80755518124807
$ws.Cells.Item(11, 6).Value = 44.05279782175224
$ws.Cells.Item(11, 12).Value = 9.992544050040527
$ws.Cells.Item(11, 13).Value = 18.85880601053765
$ws.Cells.Item(11, 14).Value = 21.10610828936523
$ws.Cells.Item(12, 2).Value = 25.31765570172891
$ws.Cells.Item(12, 3).Value = 11.72190508491989
$ws.Cells.Item(12, 4).Value = 8.083713233337066
$ws.Cells.Item(12, 5).Value = 9.812680677661994
$ws.Cells.Item(12, 6).Value = 44.16997138394793
$ws.Cells.Item(12, 12).Value = 9.991484767435695
$ws.Cells.Item(12, 13).Value = 18.88924816669324
$ws.Cells.Item(12, 14).Value = 21.09379934994226
$ws.Cells.Item(13, 2).Value = 25.28957342888547
$ws.Cells.Item(13, 3).Value = 11.69125022546705
$ws.Cells.Item(13, 4).Value = 8.08481316923824
$ws.Cells.Item(13, 5).Value = 9.811578430485463
$ws.Cells.Item(13, 6).Value = 44.14468037379757
$ws.Cells.Item(13, 12).Value = 9.991702855708466
$ws.Cells.Item(13, 13).Value = 18.88266639715574
$ws.Cells.Item(13, 14).Value = 21.09643706397652
$ws.Cells.Item(14, 2).Value = 25.19795481599297
$ws.Cells.Item(14, 3).Value = 11.59071609542
$ws.Cells.Item(14, 4).Value = 8.088427164997695
$ws.Cells.Item(14, 5).Value = 9.807977589408207
$ws.Cells.Item(14, 6).Value = 44.06241259076016
$ws.Cells.Item(14, 12).Value = 9.992452374515858
$ws.Cells.Item(14, 13).Value = 18.86129898460881
$ws.Cells.Item(14, 14).Value = 21.10508966206535
$ws.Cells.Item(15, 2).Value = 25.14181030406612
$ws.Cells.Item(15, 3).Value = 11.52870389615487
$ws.Cells.Item(15, 4).Value = 8.090661473962403
$ws.Cells.Item(15, 5).Value = 9.80576722147471
$ws.Cells.Item(15, 6).Value = 44.01218538977302
$ws.Cells.Item(15, 12).Value = 9.992940898642296
$ws.Cells.Item(15, 13).Value = 18.84828582577397
$ws.Cells.Item(15, 14).Value = 21.11042837304144
$ws.Cells.Item(16, 2).Value = 24.81977724399702
$ws.Cells.Item(16, 3).Value = 11.16670058903792
$ws.Cells.Item(16, 4).Value = 8.103781145715985
$ws.Cells.Item(16, 5).Value = 9.793027502490602
$ws.Cells.Item(16, 6).Value = 43.72693248279052
$ws.Cells.Item(16, 12).Value = 9.996191502963386
$ws.Cells.Item(16, 13).Value = 18.77487572480175
$ws.Cells.Item(16, 14).Value = 21.14161566743192
$ws.Cells.Item(17, 2).Value = 24.62206928776163
$ws.Cells.Item(17, 3).Value = 10.93868745188507
$ws.Cells.Item(17, 4).Value = 8.112110867959331
$ws.Cells.Item(17, 5).Value = 9.785146737974367
$ws.Cells.Item(17, 6).Value = 43.55430721945737
$ws.Cells.Item(17, 12).Value = 9.998588240377027
$ws.Cells.Item(17, 13).Value = 18.73089323878376
$ws.Cells.Item(17, 14).Value = 21.1612768852957
$ws.Cells.Item(18, 2).Value = 24.50834127544064
$ws.Cells.Item(18, 3).Value = 10.80535339416076
$ws.Cells.Item(18, 4).Value = 8.117005050502733
$ws.Cells.Item(18, 5).Value = 9.78058982172602
$ws.Cells.Item(18, 6).Value = 43.45591424834596
$ws.Cells.Item(18, 12).Value = 10.00011476515272
$ws.Cells.Item(18, 13).Value = 18.70598796012176
$ws.Cells.Item(18, 14).Value = 21.17277960944947
$ws.Cells.Item(19, 2).Value = 24.46983715166849
$ws.Cells.Item(19, 3).Value = 10.75983264275746
$ws.Cells.Item(19, 4).Value = 8.118679833027709
$ws.Cells.Item(19, 5).Value = 9.779042778974995
$ws.Cells.Item(19, 6).Value = 43.4227565101237
$ws.Cells.Item(19, 12).Value = 10.00065703904525
$ws.Cells.Item(19, 13).Value = 18.69762335399626
$ws.Cells.Item(19, 14).Value = 21.17670755170897
$ws.Cells.Item(20, 2).Value = 24.64311784085489
$ws.Cells.Item(20, 3).Value = 10.96318605422018
$ws.Cells.Item(20, 4).Value = 8.111213474971366
$ws.Cells.Item(20, 5).Value = 9.785988151973806
$ws.Cells.Item(20, 6).Value = 43.57259127322155
$ws.Cells.Item(20, 12).Value = 9.998317788320128
$ws.Cells.Item(20, 13).Value = 18.73553477358595
$ws.Cells.Item(20, 14).Value = 21.15916382018947
$ws.Cells.Item(21, 2).Value = 25.2248709647268
$ws.Cells.Item(21, 3).Value = 11.62033543342614
$ws.Cells.Item(21, 4).Value = 8.087361345930182
$ws.Cells.Item(21, 5).Value = 9.80903623287031
$ws.Cells.Item(21, 6).Value = 44.08654256538427
$ws.Cells.Item(21, 12).Value = 9.992226091390503
$ws.Cells.Item(21, 13).Value = 18.86755951480212
$ws.Cells.Item(21, 14).Value = 21.102540108532
$ws.Cells.Item(22, 2).Value = 25.60405759415292
$ws.Cells.Item(22, 3).Value = 12.03046398361788
$ws.Cells.Item(22, 4).Value = 8.072695496653091
$ws.Cells.Item(22, 5).Value = 9.823886700131146
$ws.Cells.Item(22, 6).Value = 44.42985948230206
$ws.Cells.Item(22, 12).Value = 9.989561781480235
$ws.Cells.Item(22, 13).Value = 18.95721733290616
$ws.Cells.Item(22, 14).Value = 21.06726674223951
$ws.Cells.Item(23, 2).Value = 25.40181547170257
$ws.Cells.Item(23, 3).Value = 11.81333606067718
$ws.Cells.Item(23, 4).Value = 8.080438274486195
$ws.Cells.Item(23, 5).Value = 9.815980104165812
$ws.Cells.Item(23, 6).Value = 44.24597359325301
$ws.Cells.Item(23, 12).Value = 9.990863323108073
$ws.Cells.Item(23, 13).Value = 18.90906279063186
$ws.Cells.Item(23, 14).Value = 21.08593393798265
$ws.Cells.Item(24, 2).Value = 24.6336019761089
$ws.Cells.Item(24, 3).Value = 10.95211723403105
$ws.Cells.Item(24, 4).Value = 8.11161885835919
$ws.Cells.Item(24, 5).Value = 9.785607830207471
$ws.Cells.Item(24, 6).Value = 43.56432239380574
$ws.Cells.Item(24, 12).Value = 9.998439596805204
$ws.Cells.Item(24, 13).Value = 18.7334351481554
$ws.Cells.Item(24, 14).Value = 21.160118516427
$ws.Cells.Item(25, 2).Value = 23.80510719635138
$ws.Cells.Item(25, 3).Value = 9.939472193788083
$ws.Cells.Item(25, 4).Value = 8.149193326718898
$ws.Cells.Item(25, 5).Value = 9.75192344518921
$ws.Cells.Item(25, 6).Value = 42.86362835255515
$ws.Cells.Item(25, 12).Value = 10.01226940659424
$ws.Cells.Item(25, 13).Value = 18.55905051383595
$ws.Cells.Item(25, 14).Value = 21.247557454848
